$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the last-changed date as an Excel serial
# date number. All data rows (2 through 37) are being bumped forward by
# one day: 45662 (2025-01-05) -> 45663 (2025-01-06).
$ws.Range("C2:C37").Value = 45663
